$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the monthly IPCA values for 2023 (row 2) that were previously missing
$ws.Range("H2").Value = 0.12
$ws.Range("I2").Value = 0.23
$ws.Range("J2").Value = 0.26
$ws.Range("K2").Value = 0.24
$ws.Range("L2").Value = 0.28

# Move the "Fonte" source link text from O3 up to O2, and remove O3 entirely
$ws.Range("O2").Value = $ws.Range("O3").Value2
$ws.Range("O3").Clear()
